# Insert a new data row at row 67 (pushes existing rows 67-185 down to 68-186)
# and populate it with a new record for Ciboulette / Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 67..185 down by one row.
$ws.Rows.Item(67).Insert()

# Populate the freshly inserted row 67 with the new observation.
$ws.Cells.Item(67, 1).Value2  = 10
$ws.Cells.Item(67, 2).Value   = "Vega Modelo de Temuco"
$ws.Cells.Item(67, 3).Value   = "La Araucanía"
$ws.Cells.Item(67, 4).Value2  = 44536
$ws.Cells.Item(67, 5).Value2  = 9
$ws.Cells.Item(67, 6).Value2  = 100112039
$ws.Cells.Item(67, 7).Value   = "Ciboulette"
$ws.Cells.Item(67, 8).Value   = "Sin especificar"
$ws.Cells.Item(67, 9).Value   = "Primera"
$ws.Cells.Item(67, 10).Value2 = 65
$ws.Cells.Item(67, 11).Value2 = 5000
$ws.Cells.Item(67, 12).Value2 = 5000
$ws.Cells.Item(67, 13).Value2 = 5000
$ws.Cells.Item(67, 14).Value  = "$/docena de atados"
$ws.Cells.Item(67, 15).Value  = "Provincia de Cautín"
$ws.Cells.Item(67, 16).Value2 = 1667
$ws.Cells.Item(67, 17).Value2 = 3
$ws.Cells.Item(67, 18).Value  = "Hortaliza"
